# Regenerate s_val data to filter save games.
# Update the computed stat columns (TB, d2S, K, IP, sum) for each game row,
# leaving the date (A) and Win (F) columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    3 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    4 = @{ B = 1.505614041169197; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    5 = @{ B = 0.7287194209349384; C = 0.05231270169004087; D = 0.1529057820181812; E = 0.4998867070740569; G = 1.433824611717217 }
    6 = @{ B = 0.1554434735375247; C = 0.3375848360084654; D = 0.1529057820181812; E = 0.4998867070740569; G = 1.145820798638228 }
    7 = @{ B = 3.182878228561681; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    8 = @{ B = 1.505614041169197; C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("B$row").Value = $rowVals.B
    $ws.Range("C$row").Value = $rowVals.C
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("E$row").Value = $rowVals.E
    $ws.Range("G$row").Value = $rowVals.G
}
